$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.386.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.51%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.467.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.67%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.13%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''529.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.94%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''131.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +0.32%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.37%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.562'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +1.54%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''2.476.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.26%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.0986'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +1.74%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -3.25%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''4.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -3.42%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.324'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.78%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''2.912.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -1.11%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''58.286.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +0.34%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''21.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.13%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''0.0000133'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -0.60%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''2.483.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.44%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''10.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -2.12%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''4.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +0.65%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''318.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.61%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  +2.42%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.00%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''65.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.81%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.402'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.28%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +1.14%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.158'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.72%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''7.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -0.03%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''175.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +4.55%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''0.0₃0740'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -0.68%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''1.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.54%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -0.60%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''6.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -1.17%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +0.14%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.12%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''17.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.55%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -4.45%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''3.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -1.84%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  +2.14%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''36.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -1.24%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.805'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +4.90%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''3.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +0.24%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''269.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -1.99%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''128.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +6.48%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''4.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -3.37%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.584'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.44%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.0931'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.36%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  +0.00%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +0.94%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''16.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.10%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''1.721.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -0.93%  '
$ws.Range("E51").Style = "Normal"
